# Update "想去人数" (F column) values on both the "展览" and "全部类型" sheets.
# Row -> (old, new) mapping taken from the diff:
#   2 -> 1996, 4 -> 123, 6 -> 16, 7 -> 1663, 9 -> 667, 11 -> 103,
#   14 -> 226, 18 -> 134, 19 -> 3862, 23 -> 359, 24 -> 715, 25 -> 508,
#   26 -> 355, 28 -> 1659, 29 -> 23, 31 -> 10

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1996
    4  = 123
    6  = 16
    7  = 1663
    9  = 667
    11 = 103
    14 = 226
    18 = 134
    19 = 3862
    23 = 359
    24 = 715
    25 = 508
    26 = 355
    28 = 1659
    29 = 23
    31 = 10
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
